$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the cached "datetimeFigureOut" footer field text that shows
#    up on the slide master and every slide layout (25/11/2022 -> 02/12/2022).
# ---------------------------------------------------------------------
function Update-DateShape($sh) {
    if ($sh.HasTextFrame) {
        $tr = $sh.TextFrame.TextRange
        if ($tr.Text -eq "25/11/2022") {
            $tr.Text = "02/12/2022"
        }
    }
}

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    Update-DateShape($master.Shapes.Item($i))
}

$layouts = $master.CustomLayouts
for ($L = 1; $L -le $layouts.Count; $L++) {
    $layout = $layouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        Update-DateShape($layout.Shapes.Item($i))
    }
}

# ---------------------------------------------------------------------
# 2) Slide 1: the credits textbox now lists a third author. Replace just
#    the trailing " e Stevan Augusto" run with the new author list; the
#    textbox has spAutoFit so PowerPoint grows its height automatically
#    to fit the now two-line caption.
# ---------------------------------------------------------------------
$slide1 = $p.Slides.Item(1)
$creditsShape = $slide1.Shapes.Item(3)
$creditsRange = $creditsShape.TextFrame.TextRange

$oldTail = " e Stevan Augusto"
$newTail = ", Stevan Augusto e Henrique Machado"
$startPos = $creditsRange.Text.IndexOf($oldTail) + 1

$creditsRange.Characters($startPos, $oldTail.Length).Text = $newTail
